$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The ferrite_length parameter value (row 7, column B) is updated from "109 mm" to "218 mm"
$ws.Range("B7").Value = "218 mm"

# Reflect the new active cell selection recorded in the saved file
$ws.Range("B8").Select()
